$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.178.61"
$ws.Range("E2").Value = "  +3.08%  "
$ws.Range("D3").Value = "2.307.94"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'310.56"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").Value = "'101.40"
$ws.Range("E6").Value = "  +6.28%  "
$ws.Range("D7").Value = "'0.537"
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.525"
$ws.Range("E9").Value = "  +6.93%  "
$ws.Range("D10").Value = "'35.99"
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("E11").Value = "  +3.59%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "'7.16"
$ws.Range("E13").Value = "  +7.34%  "
$ws.Range("D14").Value = "2.664.18"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("D16").Value = "2.305.20"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "'0.811"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "43.085.92"
$ws.Range("E18").Value = "  +3.06%  "
$ws.Range("D19").Value = "'12.56"
$ws.Range("E19").Value = "  +1.63%  "
$ws.Range("D20").Value = "0.0₃0922"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").Value = "'68.61"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").Value = "'241.20"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  +4.10%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.63"
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'24.83"
$ws.Range("E27").Value = "  +5.00%  "
$ws.Range("D28").Value = "'37.63"
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'9.66"
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.11"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "'167.62"
$ws.Range("E31").Value = "  +4.59%  "
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").Value = "'17.73"
$ws.Range("E35").Value = "  +4.54%  "
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("E37").Value = "  +3.21%  "
$ws.Range("D38").Value = "'2.40"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").Value = "'1.83"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").Value = "'4.31"
$ws.Range("E41").Value = "  +8.02%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "1.985.65"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("D45").Value = "'19.19"
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("E46").Value = "  +3.71%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "'2.96"
$ws.Range("E48").Value = "  +18.12%  "
$ws.Range("D49").Value = "'55.69"
$ws.Range("E49").Value = "  +5.23%  "
$ws.Range("D50").Value = "2.531.53"
$ws.Range("E51").Value = "  +2.18%  "
